# Auto-generated edit script: updates crypto price/volume table cells
# per the commit "Updated cryptos list on Sat Jul 20 22:34:18 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.139.78'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '3.518.74'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.32'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.02'
$ws.Range('E6').Value = '  +3.00%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.596'
$ws.Range('E8').Value = '  +3.58%  '
$ws.Range('E9').Value = '  +7.11%  '
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '4.128.53'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.04'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '67.133.06'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('D17').Value = '3.522.72'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.29'
$ws.Range('E19').Value = '  +1.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '397.25'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000122'
$ws.Range('E25').Value = '  -3.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.22'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.29'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.41'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.64'
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '163.70'
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.894'
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.91'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('E38').Value = '  +3.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.83'
$ws.Range('E39').Value = '  +5.67%  '
$ws.Range('E40').Value = '  +1.16%  '
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.49'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.803.60'
$ws.Range('E43').Value = '  -1.28%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.61'
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.86'
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('E46').Value = '  -2.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '339.95'
$ws.Range('E47').Value = '  -3.34%  '
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.54'
$ws.Range('E50').Value = '  +0.40%  '
$ws.Range('E51').Value = '  -1.01%  '
